$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 112 and 113: row 112 gets old row 113 values, row 113 gets old row 112 values
$ws.Range("B112").Value = 6779681
$ws.Range("C112").Value = "Denmark Superligaen"
$ws.Range("D112").Value = "Denmark Superligaen"
$ws.Range("E112").Value = 45340.41666666666
$ws.Range("F112").Value = "FC Nordsjaelland"
$ws.Range("G112").Value = "Lyngby"
$ws.Range("H112").Value = 3
$ws.Range("I112").Value = 2
$ws.Range("J112").Value = "H"
$ws.Range("K112").Value = 1.533
$ws.Range("L112").Value = 4.333
$ws.Range("M112").Value = 5
$ws.Range("N112").Value = 1.4
$ws.Range("O112").Value = 5.5
$ws.Range("P112").Value = 6.5
$ws.Range("Q112").Value = -1.25
$ws.Range("R112").Value = 1.85
$ws.Range("S112").Value = 2
$ws.Range("T112").Value = 3
$ws.Range("U112").Value = 1.95
$ws.Range("V112").Value = 1.9
$ws.Range("W112").Value = 0.3999999999999999
$ws.Range("X112").Value = -1
$ws.Range("Y112").Value = -1
$ws.Range("Z112").Value = -0.5
$ws.Range("AA112").Value = 0.5
$ws.Range("AB112").Value = 0.95
$ws.Range("AC112").Value = -1

$ws.Range("B113").Value = 6779678
$ws.Range("C113").Value = "Denmark Superligaen"
$ws.Range("D113").Value = "Denmark Superligaen"
$ws.Range("E113").Value = 45340.41666666666
$ws.Range("F113").Value = "Hvidovre IF"
$ws.Range("G113").Value = "Randers FC"
$ws.Range("H113").Value = 1
$ws.Range("I113").Value = 3
$ws.Range("J113").Value = "A"
$ws.Range("K113").Value = 3.3
$ws.Range("L113").Value = 3.6
$ws.Range("M113").Value = 2.05
$ws.Range("N113").Value = 4.5
$ws.Range("O113").Value = 3.8
$ws.Range("P113").Value = 1.8
$ws.Range("Q113").Value = 0.75
$ws.Range("R113").Value = 1.85
$ws.Range("S113").Value = 2
$ws.Range("T113").Value = 2.5
$ws.Range("U113").Value = 1.9
$ws.Range("V113").Value = 1.95
$ws.Range("W113").Value = -1
$ws.Range("X113").Value = -1
$ws.Range("Y113").Value = 0.8
$ws.Range("Z113").Value = -1
$ws.Range("AA113").Value = 1
$ws.Range("AB113").Value = 0.8999999999999999
$ws.Range("AC113").Value = -1

# Swap row 118 and 119: row 118 gets old row 119 values, row 119 gets old row 118 values
$ws.Range("B118").Value = 6779686
$ws.Range("C118").Value = "Denmark Superligaen"
$ws.Range("D118").Value = "Denmark Superligaen"
$ws.Range("E118").Value = 45347.41666666666
$ws.Range("F118").Value = "Vejle"
$ws.Range("G118").Value = "Silkeborg IF"
$ws.Range("H118").Value = 2
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = "H"
$ws.Range("K118").Value = 3
$ws.Range("L118").Value = 3.4
$ws.Range("M118").Value = 2.3
$ws.Range("N118").Value = 3
$ws.Range("O118").Value = 3.4
$ws.Range("P118").Value = 2.3
$ws.Range("Q118").Value = 0.25
$ws.Range("R118").Value = 1.85
$ws.Range("S118").Value = 2
$ws.Range("T118").Value = 2.25
$ws.Range("U118").Value = 1.825
$ws.Range("V118").Value = 2.025
$ws.Range("W118").Value = 2
$ws.Range("X118").Value = -1
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = 0.8500000000000001
$ws.Range("AA118").Value = -1
$ws.Range("AB118").Value = -0.5
$ws.Range("AC118").Value = 0.5125

$ws.Range("B119").Value = 6780974
$ws.Range("C119").Value = "Denmark Superligaen"
$ws.Range("D119").Value = "Denmark Superligaen"
$ws.Range("E119").Value = 45347.41666666666
$ws.Range("F119").Value = "Hvidovre IF"
$ws.Range("G119").Value = "Viborg"
$ws.Range("H119").Value = 2
$ws.Range("I119").Value = 2
$ws.Range("J119").Value = "D"
$ws.Range("K119").Value = 4
$ws.Range("L119").Value = 3.7
$ws.Range("M119").Value = 1.833
$ws.Range("N119").Value = 5.25
$ws.Range("O119").Value = 4.2
$ws.Range("P119").Value = 1.615
$ws.Range("Q119").Value = 1
$ws.Range("R119").Value = 1.8
$ws.Range("S119").Value = 2.05
$ws.Range("T119").Value = 2.75
$ws.Range("U119").Value = 2
$ws.Range("V119").Value = 1.85
$ws.Range("W119").Value = -1
$ws.Range("X119").Value = 3.2
$ws.Range("Y119").Value = -1
$ws.Range("Z119").Value = 0.8
$ws.Range("AA119").Value = -1
$ws.Range("AB119").Value = 1
$ws.Range("AC119").Value = -1

# Swap row 137 and 140: row 137 gets old row 140 values, row 140 gets old row 137 values
$ws.Range("B137").Value = 6779698
$ws.Range("C137").Value = "Denmark Superligaen"
$ws.Range("D137").Value = "Denmark Superligaen"
$ws.Range("E137").Value = 45368.54166666666
$ws.Range("F137").Value = "Brondby"
$ws.Range("G137").Value = "Silkeborg IF"
$ws.Range("H137").Value = 4
$ws.Range("I137").Value = 1
$ws.Range("J137").Value = "H"
$ws.Range("K137").Value = 1.615
$ws.Range("L137").Value = 3.8
$ws.Range("M137").Value = 5.25
$ws.Range("N137").Value = 1.615
$ws.Range("O137").Value = 3.8
$ws.Range("P137").Value = 5.75
$ws.Range("Q137").Value = -0.75
$ws.Range("R137").Value = 1.85
$ws.Range("S137").Value = 2.05
$ws.Range("T137").Value = 2.5
$ws.Range("U137").Value = 1.85
$ws.Range("V137").Value = 2
$ws.Range("W137").Value = 0.615
$ws.Range("X137").Value = -1
$ws.Range("Y137").Value = -1
$ws.Range("Z137").Value = 0.8500000000000001
$ws.Range("AA137").Value = -1
$ws.Range("AB137").Value = 0.8500000000000001
$ws.Range("AC137").Value = -1

$ws.Range("B140").Value = 6779697
$ws.Range("C140").Value = "Denmark Superligaen"
$ws.Range("D140").Value = "Denmark Superligaen"
$ws.Range("E140").Value = 45368.54166666666
$ws.Range("F140").Value = "AGF Aarhus"
$ws.Range("G140").Value = "Hvidovre IF"
$ws.Range("H140").Value = 1
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = "H"
$ws.Range("K140").Value = 1.363
$ws.Range("L140").Value = 4.75
$ws.Range("M140").Value = 7.5
$ws.Range("N140").Value = 1.533
$ws.Range("O140").Value = 4.2
$ws.Range("P140").Value = 6
$ws.Range("Q140").Value = -1
$ws.Range("R140").Value = 1.875
$ws.Range("S140").Value = 1.975
$ws.Range("T140").Value = 2.5
$ws.Range("U140").Value = 1.925
$ws.Range("V140").Value = 1.925
$ws.Range("W140").Value = 0.5329999999999999
$ws.Range("X140").Value = -1
$ws.Range("Y140").Value = -1
$ws.Range("Z140").Value = 0
$ws.Range("AA140").Value = -0
$ws.Range("AB140").Value = -1
$ws.Range("AC140").Value = 0.925

# Row 147 updates
$ws.Range("N147").Value = 3.1
$ws.Range("P147").Value = 2.25
$ws.Range("R147").Value = 1.89
$ws.Range("S147").Value = 2.01
$ws.Range("U147").Value = 1.875
$ws.Range("V147").Value = 1.975

# Row 148 updates
$ws.Range("N148").Value = 2.9
$ws.Range("R148").Value = 2.09
$ws.Range("S148").Value = 1.81

# Row 149 updates
$ws.Range("N149").Value = 4.75
$ws.Range("O149").Value = 3.75
$ws.Range("P149").Value = 1.75
$ws.Range("R149").Value = 1.88
$ws.Range("S149").Value = 2.02
$ws.Range("U149").Value = 1.925
$ws.Range("V149").Value = 1.925

# Row 150 updates
$ws.Range("O150").Value = 3.5
$ws.Range("P150").Value = 2.4
$ws.Range("R150").Value = 1.82
$ws.Range("S150").Value = 2.08
$ws.Range("U150").Value = 1.95
$ws.Range("V150").Value = 1.9

# --- New Row 151 ---
$ws.Range("A150").Copy()
$ws.Range("A151").PasteSpecial(-4122)
$ws.Range("A151").Value = 149
$ws.Range("B151").Value = 7984001
$ws.Range("C151").Value = "Denmark Superligaen"
$ws.Range("D151").Value = "Denmark Superligaen"
$ws.Range("E150").Copy()
$ws.Range("E151").PasteSpecial(-4122)
$ws.Range("E151").Value = 45389.54166666666
$ws.Range("F151").Value = "AGF Aarhus"
$ws.Range("G151").Value = "Midtjylland"
$ws.Range("K151").Value = 3.1
$ws.Range("L151").Value = 3.4
$ws.Range("M151").Value = 2.25
$ws.Range("N151").Value = 3.25
$ws.Range("O151").Value = 3.4
$ws.Range("P151").Value = 2.2
$ws.Range("Q151").Value = 0.25
$ws.Range("R151").Value = 1.98
$ws.Range("S151").Value = 1.92
$ws.Range("T151").Value = 2.25
$ws.Range("U151").Value = 1.85
$ws.Range("V151").Value = 2
$ws.Range("W151").Value = 0
$ws.Range("X151").Value = 0
$ws.Range("Y151").Value = 0
$ws.Range("Z151").Value = 0
$ws.Range("AA151").Value = 0

# --- New Row 152 ---
$ws.Range("A150").Copy()
$ws.Range("A152").PasteSpecial(-4122)
$ws.Range("A152").Value = 150
$ws.Range("B152").Value = 7984002
$ws.Range("C152").Value = "Denmark Superligaen"
$ws.Range("D152").Value = "Denmark Superligaen"
$ws.Range("E150").Copy()
$ws.Range("E152").PasteSpecial(-4122)
$ws.Range("E152").Value = 45390.58333333334
$ws.Range("F152").Value = "Brondby"
$ws.Range("G152").Value = "Silkeborg IF"
$ws.Range("K152").Value = 1.5
$ws.Range("L152").Value = 4.333
$ws.Range("M152").Value = 6.5
$ws.Range("N152").Value = 1.444
$ws.Range("O152").Value = 4.5
$ws.Range("P152").Value = 7.5
$ws.Range("Q152").Value = -1.25
$ws.Range("R152").Value = 2.08
$ws.Range("S152").Value = 1.82
$ws.Range("T152").Value = 2.75
$ws.Range("U152").Value = 1.9
$ws.Range("V152").Value = 1.95
$ws.Range("W152").Value = 0
$ws.Range("X152").Value = 0
$ws.Range("Y152").Value = 0
$ws.Range("Z152").Value = 0
$ws.Range("AA152").Value = 0
